$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '49.888.68'
$ws.Range('E2').Value = '  +3.83%  '
$ws.Range('D3').Value = '2.646.10'
$ws.Range('E3').Value = '  +6.04%  '
$origStyle = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = $origStyle
$ws.Range('E4').Value = '  +0.02%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '114.09'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +7.92%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '326.34'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  +2.06%  '
$ws.Range('E7').Value = '  +1.61%  '
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  +0.02%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.554'
$ws.Range('D9').Style = $origStyle
$ws.Range('E9').Value = '  +3.09%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.13'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  +6.27%  '
$ws.Range('E11').Value = '  +0.28%  '
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0821'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  +2.17%  '
$ws.Range('E13').Value = '  +1.00%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.37'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  +3.98%  '
$ws.Range('D15').Value = '3.060.74'
$ws.Range('E15').Value = '  +5.97%  '
$ws.Range('D16').Value = '2.645.72'
$ws.Range('E16').Value = '  +5.55%  '
$ws.Range('E17').Value = '  +4.54%  '
$ws.Range('D18').Value = '49.791.71'
$ws.Range('E18').Value = '  +3.96%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.18'
$ws.Range('D19').Style = $origStyle
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('E21').Value = '  +2.04%  '
$ws.Range('E22').Value = '  +2.62%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.06'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  +1.29%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '277.36'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  +1.93%  '
$ws.Range('E25').Value = '  +2.77%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.79'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  +4.06%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.99'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  +2.82%  '
$ws.Range('E29').Value = '  -2.98%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.14'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  +4.15%  '
$ws.Range('E31').Value = '  +0.44%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.34'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +2.46%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.43'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  +2.57%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.56'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +2.36%  '
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('E37').Value = '  +7.21%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.88'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +7.13%  '
$ws.Range('E39').Value = '  +8.24%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '126.10'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +3.08%  '
$ws.Range('E41').Value = '  +1.79%  '
$ws.Range('E42').Value = '  +1.59%  '
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.08'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  -0.50%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0315'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  +3.56%  '
$ws.Range('D45').Value = '2.080.41'
$ws.Range('E45').Value = '  +4.14%  '
$ws.Range('E46').Value = '  +5.72%  '
$ws.Range('E47').Value = '  +14.17%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.98'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +4.71%  '
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.11'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  +2.46%  '
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.38'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  +4.03%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '60.31'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  +7.48%  '
